# "Them tt chuong sach TA" - add chapter title text for the English-language
# textbook reference on the title slide.
#
# Slide 1, shape "Rectangle 3" (subTitle placeholder) holds two paragraphs:
#   1) "Chương 5.Mô hình nhị phân độc lập"
#   2) "IIR.Chap11.Probabilistic information retrieval"
#
# Paragraph 2 is updated so the short-form chapter reference "Chap11"
# becomes "C11" -> "IIR.C11.Probabilistic information retrieval".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 3")
$tr = $shp.TextFrame.TextRange

# Paragraph 1: re-type the " 5.Mô " portion (originally two adjacent runs:
# " " and "5.Mô ") as a single run so it merges into one text run.
$para1 = $tr.Paragraphs(1)
$seg1 = $para1.Characters(7, 6)
$seg1.Text = " 5.Mô "

# Paragraph 2: change "IIR.Chap11.Probabilistic " (chars 1-25) so the run
# reads "IIR.C11.Probabilistic " and splits away from the unchanged
# "information retrieval" tail.
$para2 = $tr.Paragraphs(2)
$seg2 = $para2.Characters(1, 25)
$seg2.Text = "IIR.C11.Probabilistic "
